$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.151.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.945.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4816"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.63%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2917"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06807"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "104.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.34%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.971.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.33%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07863"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.309"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6900"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.163.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.216.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007621"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.593"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.4743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "

# Row 24
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.466"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.591"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.45%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.68%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.45%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.138"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1016"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.09%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.391"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.64%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.636"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.83%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.540"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.56%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.361"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04841"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7433"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.38%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.33%  "

# Row 38
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.730"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.33%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.602"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.49%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.649"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "77.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.037"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.56%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8724"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4378"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.64%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.028.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.33%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.583"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.93%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.239"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1214"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
